$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 0.135186667482864
$ws.Range("C7").Value = 0.9840273120519328
$ws.Range("D7").Value = 5.006397119580056
$ws.Range("E7").Value = 2.237497959681764
$ws.Range("F7").Value = 2.263390329593289
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.1814798078799119
$ws.Range("C8").Value = 1.020846580860502
$ws.Range("D8").Value = 5.181999496594389
$ws.Range("E8").Value = 2.276400557150343
$ws.Range("F8").Value = 2.300455204406243
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = 0.09581464916219476
$ws.Range("C9").Value = 1.411823298613351
$ws.Range("D9").Value = 8.906796727195324
$ws.Range("E9").Value = 2.984425694701633
$ws.Range("F9").Value = 3.060377733858364
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.6619987679365094
$ws.Range("C10").Value = 1.043113688983639
$ws.Range("D10").Value = 5.226911917360573
$ws.Range("E10").Value = 2.286244063384435
$ws.Range("F10").Value = 2.277657863755817
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.1186778033967304
$ws.Range("C11").Value = 0.3680761991877547
$ws.Range("D11").Value = 0.2032499457136513
$ws.Range("E11").Value = 0.4508325029472158
$ws.Range("F11").Value = 0.4862683475903228
$ws.Range("G11").Value = 5
